$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.706.80"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "2.116.71"
$ws.Range("E3").Value = "  +10.10%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "255.41"
$ws.Range("E5").Value = "  +2.28%  "

$ws.Range("D6").Value = "0.667"
$ws.Range("E6").Value = "  -4.36%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "47.19"
$ws.Range("E8").Value = "  +6.32%  "

$ws.Range("D9").Value = "61.49"
$ws.Range("E9").Value = "  +4.91%  "

$ws.Range("D10").Value = "0.376"
$ws.Range("E10").Value = "  +2.17%  "

$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").Value = "2.428.55"
$ws.Range("E13").Value = "  +10.36%  "

$ws.Range("D14").Value = "14.37"
$ws.Range("E14").Value = "  -1.75%  "

$ws.Range("D15").Value = "0.836"
$ws.Range("E15").Value = "  +4.32%  "

$ws.Range("D16").Value = "2.123.86"
$ws.Range("E16").Value = "  +10.51%  "

$ws.Range("D17").Value = "5.14"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "36.754.98"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").Value = "73.86"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").Value = "13.28"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "241.48"
$ws.Range("E22").Value = "  -4.12%  "

$ws.Range("D23").Value = "5.21"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -6.83%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "172.68"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.77"
$ws.Range("E27").Value = "  +15.50%  "

$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  +4.64%  "

$ws.Range("D29").Value = "2.03"
$ws.Range("E29").Value = "  -7.86%  "

$ws.Range("D30").Value = "28.84"
$ws.Range("E30").Value = "  +62.48%  "

$ws.Range("E31").Value = "  -4.42%  "

$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "0.0955"
$ws.Range("E33").Value = "  +10.12%  "

$ws.Range("E34").Value = "  -2.48%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.965"
$ws.Range("E35").Value = "  +8.11%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").Value = "  +16.96%  "

$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -5.07%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("E39").Value = "  -4.58%  "

$ws.Range("E40").Value = "  -10.99%  "

$ws.Range("E41").Value = "  +7.77%  "

$ws.Range("D42").Value = "0.0226"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").Value = "99.17"
$ws.Range("E43").Value = "  -6.40%  "

$ws.Range("D44").Value = "2.79"
$ws.Range("E44").Value = "  +8.90%  "

$ws.Range("D45").Value = "16.22"
$ws.Range("E45").Value = "  -6.68%  "

$ws.Range("D46").Value = "1.358.11"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").Value = "7.26"
$ws.Range("E47").Value = "  +12.66%  "

$ws.Range("D48").Value = "0.0842"
$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("D49").Value = "2.313.65"
$ws.Range("E49").Value = "  +10.23%  "

$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  +1.48%  "
